$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended after the last existing row (row 191), matching the
# "Actualizacion desde MV -datos-" update (05, 06 and 07-10-2021).
$newRows = @(
    @("05-10-2021", 4131, 549, 501, 933, 1013, 1135),
    @("06-10-2021", 5855, 1130, 764, 1186, 1347, 1427),
    @("07-10-2021", 4088, 565, 723, 1176, 778, 846)
)

$startRow = 192
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a date-shaped label ("dd-mm-yyyy"). Writing it straight
    # to .Value would make Excel auto-detect it as a real date and stamp the
    # cell with a date NumberFormat/style. Instead build it as a text formula
    # and immediately convert it to its literal value (copy / paste-special
    # values), which is how the source data keeps these as plain shared
    # strings with the default (unstyled) cell format.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = '="' + $row[0] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
